$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the G2:G3 merged cell - G2 and G3 become independent cells
$ws.Range("G2:G3").UnMerge()

# Row 2 header text changes
$ws.Range("G2").Value = "第二级"
$ws.Range("H2").Value = "第二级"

# Row 3 text changes (G3/H3 first, C3 later to control shared-string ordering)
$ws.Range("G3").Value = "第三级"
$ws.Range("H3").Value = "第三级"

# New data rows
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "name10"
$ws.Range("C4").Value = $True
$ws.Range("D4").Value = 13
$ws.Range("E4").Value = 14
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 17

$ws.Range("A5").Value = 20
$ws.Range("B5").Value = "name20"
$ws.Range("C5").Value = $False
$ws.Range("D5").Value = 23
$ws.Range("E5").Value = 24
$ws.Range("F5").Value = 25
$ws.Range("G5").Value = 26
$ws.Range("H5").Value = 27

$ws.Range("C3").Value = "空 格"

# Selection moves to C3
[void]$ws.Range("C3").Select()
